$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "#SOS3"
$ws.Range("A13").Value = "#EOS3"
$ws.Range("A12").Value = "게임을 종료합니다."

$ws.Range("A12").Select()
